{"js": "// Load all paragraphs of the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// 1) Title paragraph: \"Logs de \" + \"t\u00e2ches\" (two runs, spell-check markers)\n//    -> single run \"Logs de t\u00e2ches\" with the proofErr markers gone.\n//    A plain text replace on this paragraph leaves a stray trailing\n//    <w:proofErr/> behind (it sits after the last run), so instead we\n//    insert a brand-new paragraph with the final text/formatting right\n//    before it and delete the old paragraph outright.\nconst titlePara = paragraphs.items[4];\ntitlePara.insertParagraph(\"Logs de t\u00e2ches\", Word.InsertLocation.before);\ntitlePara.delete();\nawait context.sync();\n\n// Re-resolve the paragraph collection; indices are stable across the\n// delete/insert pair above (one paragraph removed, one added in its place).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\n// 2) \"Ajout des classes de base, cr\u00e9ation de namespace (...)\" paragraph:\n//    three runs (\"... de \", \"namespace\", \" (1 heure...)\") collapse into one.\n//    The proofErr pair here is fully interior (followed by more run\n//    content), so a simple in-place text replace cleans it up completely.\nconst namespacePara = paras.items[7];\nnamespacePara.insertText(\n  \"Ajout des classes de base, cr\u00e9ation de namespace (1 heure, TBH, CL, 28 octobre)\",\n  Word.InsertLocation.replace\n);\n\n// 3) \"Singleton template et non copyable (...)\" paragraph: five runs\n//    collapse into one, same proofErr situation as above (fully interior).\nconst singletonPara = paras.items[8];\nsingletonPara.insertText(\n  \"Singleton template et non copyable (10 minutes, CL, 29 octobre)\",\n  Word.InsertLocation.replace\n);\n\n// 4) New bullet after \"Classe projectile (...)\": add a sibling list item\n//    \"D\u00e9lai de tir du joueur (TBH, 10 minutes)\". insertParagraph copies the\n//    paragraph style/numbering (Paragraphedeliste, numId 2, jc both) from\n//    the paragraph it is called on.\nconst projectilePara = paras.items[9];\nprojectilePara.insertParagraph(\n  \"D\u00e9lai de tir du joueur (TBH, 10 minutes)\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Journal de projet \u2014 \"Debut de fonction tir dans player\"\n#\n# 1) Title paragraph \"Logs de \" + \"t\u00e2ches\" (two runs + spellcheck markers)\n#    -> single run \"Logs de t\u00e2ches\".\n# 2) \"Ajout des classes de base, cr\u00e9ation de namespace (...)\" paragraph:\n#    three runs -> one run.\n# 3) \"Singleton template et non copyable (...)\" paragraph: five runs -> one run.\n# 4) New bullet \"D\u00e9lai de tir du joueur (TBH, 10 minutes)\" right after\n#    \"Classe projectile (...)\".\n\n$d = $word.ActiveDocument\n\n# --- 1) Title paragraph --------------------------------------------------\n# A plain in-place text replace merges the runs but this engine leaves a\n# stray trailing <w:proofErr/> behind (it is the last child of the\n# paragraph, after the last run). Side-stepping that: splice in a brand\n# new paragraph carrying the final text/formatting right before the old\n# one, then delete the old paragraph (proofErr and all) outright.\n$titlePara = $d.Paragraphs.Item(5)\n$titlePara.Range.InsertParagraphBefore()\n$newTitlePara = $d.Paragraphs.Item(5)\n$newTitlePara.Range.Text = \"Logs de t\u00e2ches\"\n$oldTitlePara = $d.Paragraphs.Item(6)\n$oldTitlePara.Range.Delete()\n\n# --- 2) \"Ajout des classes de base...\" paragraph -------------------------\n# Here the proofErr pair sits entirely between runs (more run content\n# follows), so an in-place replace cleans it up fully. Because the\n# resulting visible text is identical to the original, a direct\n# single-shot replace is treated as a no-op by this engine, so first push\n# through a distinct placeholder, then overwrite it with the real text \u2014\n# this guarantees the run-merge/proofErr-cleanup actually happens.\n$nsPara = $d.Paragraphs.Item(8)\n$nsRange = $nsPara.Range\n$nsRange.MoveEnd(1, -1)\n$nsRange.Text = \"~~tmp~~\"\n$nsRange2 = $nsPara.Range\n$nsRange2.MoveEnd(1, -1)\n$nsRange2.Text = \"Ajout des classes de base, cr\u00e9ation de namespace (1 heure, TBH, CL, 28 octobre)\"\n\n# --- 3) \"Singleton template et non copyable...\" paragraph ----------------\n$stPara = $d.Paragraphs.Item(9)\n$stRange = $stPara.Range\n$stRange.MoveEnd(1, -1)\n$stRange.Text = \"~~tmp~~\"\n$stRange2 = $stPara.Range\n$stRange2.MoveEnd(1, -1)\n$stRange2.Text = \"Singleton template et non copyable (10 minutes, CL, 29 octobre)\"\n\n# --- 4) New bullet right after \"Classe projectile (...)\" -----------------\n$projPara = $d.Paragraphs.Item(10)\n$projPara.Range.InsertParagraphAfter()\n$newBulletPara = $d.Paragraphs.Item(11)\n$newBulletRange = $newBulletPara.Range\n$newBulletRange.MoveEnd(1, -1)\n$newBulletRange.Text = \"D\u00e9lai de tir du joueur (TBH, 10 minutes)\"\n"}
